# "Addition of demo tour"
# Adds two new worksheets ("createNewAccount" and "Registration") that
# mirror the Flag/Yes header pattern used on the existing sheets, and
# tweaks the active-sheet / selection state left behind in the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # CreateAccount
$ws2 = $wb.Worksheets.Item(2)   # ContactUs

# --- New sheet: createNewAccount -------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "createNewAccount"

$ws3.Range("A1").Value = "Flag"
$ws3.Range("A2").Value = "Yes"

# Highlight the header cell in yellow
$ws3.Range("A1").Interior.Color = 65535

# --- New sheet: Registration -------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Registration"

$ws4.Range("A1").Value = "Flag"
$ws4.Range("B1").Value = "TestCaseName"
$ws4.Range("A2").Value = "Yes"
$ws4.Range("B2").Value = "TC01"

# Reuse the existing header/value formatting from CreateAccount
[void]$ws1.Range("A1").Copy()
[void]$ws4.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats

[void]$ws1.Range("A2").Copy()
[void]$ws4.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$ws4.Columns.Item(2).ColumnWidth = 14.140625
[void]$ws4.Range("D9").Select()

# --- Leftover UI state on the pre-existing sheets -----------------------
[void]$ws2.Range("A1:C2").Select()
[void]$ws1.Range("E24").Select()
